# Add three new form-submission rows to the raw collection-results sheet
# ("八位序列号收集收集结果yd5"), mirroring three new responses that came in
# after the existing 104 rows of data (submitter, timestamp, serial number,
# QQ number).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("八位序列号收集收集结果yd5")

$dateFormat = "yyyy/m/d h:mm:ss;@"

$rows = @(
    @{ Row = 105; A = "刻";      B = 45973.8885185185; C = "ba3ef151"; D = "2512200699" },
    @{ Row = 106; A = "⁦⁦";      B = 45973.9549305556; C = "5525cf80"; D = "2628639578" },
    @{ Row = 107; A = "鱼丸粗面"; B = 45974.0308449074; C = "d2051fdf"; D = "2942666153" }
)

foreach ($r in $rows) {
    $i = $r.Row

    # Column A - submitter (plain text, not purely numeric so it stores as text naturally)
    $ws.Cells.Item($i, 1).Value = $r.A

    # Column B - submission timestamp, stored as a date/time serial number
    $ws.Cells.Item($i, 2).Value = $r.B
    $ws.Cells.Item($i, 2).NumberFormat = $dateFormat

    # Columns C/D - serial number and QQ number. These are digit-heavy values
    # that must stay text (matching the rest of the sheet), so build them via
    # a text formula and then paste back as plain values - this avoids Excel
    # re-typing them as numbers or tagging the cell with a "number stored as
    # text" quote-prefix style.
    $ws.Cells.Item($i, 3).Formula = '="' + $r.C + '"'
    $ws.Cells.Item($i, 4).Formula = '="' + $r.D + '"'
}

$dataRange = $ws.Range("C105:D107")
$dataRange.Copy()
$dataRange.PasteSpecial(-4163)
